$p = $ppt.ActivePresentation

# The deck originally had 23 slides. This edit trims it down to the
# slides that were actually presented:
#   - slide 6  ("Q & A") is removed
#   - slides 15-23 (extra "largest/smallest set ..." topic slides that
#     go beyond "Minimum edge dominating set") are removed
# All other slides keep their relative order.
#
# Delete from the highest index to the lowest so that the indices of
# slides not yet processed stay valid as we go.
$positionsToDelete = @(23, 22, 21, 20, 19, 18, 17, 16, 15, 6)

foreach ($pos in $positionsToDelete) {
    $p.Slides.Item($pos).Delete()
}
